$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "ser: N" blog-card entries shift: 164 is dropped, 165/167 slide down
# one slot, and a new "ser: 168" entry is appended. On the sheet this shows
# up as the three blog cells in row 10 (H10, D10, B10) advancing to the
# next ser number in the sequence.
$ws.Range("H10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 165"
$ws.Range("D10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 167"
$ws.Range("B10").Value = "type: blog`nwidth: 2`nheight: 1`nser: 168"

# Move the active selection from C10 to H10.
$ws.Range("H10").Select() | Out-Null
